# Commit "Add files via upload": workbook re-saved from a newer Excel build.
# The only content-level change versus the previous upload is the worksheet
# tab name: "porcatroia" -> "dataset". (Everything else in the diff -
# fileVersion/rupBuild, xr:* revision namespaces, absPath, column-width
# rounding, x14ac:dyDescent hints - is Excel-version metadata noise from the
# resave, not a data/formatting edit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "dataset"
